$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.886.81'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '1.700.25'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  -0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.95'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4021'
$ws.Range("E7").Value = '  +2.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4064'
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("E9").Value = '  -0.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.83'
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.465'
$ws.Range("E11").Value = '  -3.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08820'
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.79'
$ws.Range("E13").Value = '  +5.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.480'
$ws.Range("E14").Value = '  -1.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.045'
$ws.Range("E15").Value = '  +0.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001345'
$ws.Range("E16").Value = '  -1.43%  '
$ws.Range("D17").Value = '1.764.59'
$ws.Range("E17").Value = '  +3.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '96.69'
$ws.Range("E18").Value = '  -1.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07195'
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.94'
$ws.Range("E20").Value = '  +5.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.239'
$ws.Range("E21").Value = '  -2.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.004'
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.56'
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("D24").Value = '24.893.51'
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.892'
$ws.Range("E26").Value = '  -5.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.655'
$ws.Range("E27").Value = '  +27.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.11'
$ws.Range("E28").Value = '  +1.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.41'
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '143.72'
$ws.Range("E30").Value = '  +4.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.173'
$ws.Range("E31").Value = '  -3.77%  '
$ws.Range("D32").Value = '1.986.33'
$ws.Range("E32").Value = '  +5.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.280'
$ws.Range("E33").Value = '  +14.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08741'
$ws.Range("E34").Value = '  -1.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.409'
$ws.Range("E35").Value = '  -0.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.03169'
$ws.Range("E36").Value = '  +8.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.038'
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2847'
$ws.Range("E38").Value = '  +4.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8494'
$ws.Range("E39").Value = '  +7.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.87'
$ws.Range("E40").Value = '  +0.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09437'
$ws.Range("E41").Value = '  +3.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.05'
$ws.Range("E42").Value = '  -2.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.471'
$ws.Range("E43").Value = '  +0.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.75'
$ws.Range("E44").Value = '  +6.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.710'
$ws.Range("E45").Value = '  +5.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7460'
$ws.Range("E46").Value = '  +3.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.237'
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.405'
$ws.Range("E48").Value = '  +4.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.006'
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.36'
$ws.Range("E50").Value = '  +1.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08354'
$ws.Range("E51").Value = '  +4.56%  '
